$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "09/20/2025"
$ws.Range("B19").Value = 0.1258563147191275
$ws.Range("C19").Value = 0.8741436852808725
